$d = $word.ActiveDocument

# 1) "Lead analyst..., Pierce Protocols" bullet about the Filemaker database:
#    drop the trailing ", and creating documentation" clause.
$found1 = $d.Content.Find.Execute(
    "writing Filemaker scripts to perform complex queries and create dynamic reports, and creating documentation.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "writing Filemaker scripts to perform complex queries and create dynamic reports.",
    2)

# 2) "Mathematics researcher..." bullet about leading classes:
#    drop the trailing ", and discussing weekly tasks" clause.
$found2 = $d.Content.Find.Execute(
    "marking and giving feedback on assessed work (providing 1-1 support), and discussing weekly tasks.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "marking and giving feedback on assessed work (providing 1-1 support).",
    2)

# 3) "Junior Geometry Seminar Organiser" entry description - reword.
$found3 = $d.Content.Find.Execute(
    "Founded seminars, organised speakers, booked rooms, directed discussion and published notes.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Founded seminar, organised speakers, logistics, and directed discussion.",
    2)

# 4) Remove the entire "IT Officer, Durham University Canoe Club" entry: its
#    "2008 - 2009" date-term paragraph and the definition paragraph right after
#    it (which itself holds three runs/lines about the IT Officer role).
$paras = $d.Paragraphs
$defIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Text -like "*IT Officer, Durham University Canoe Club*") {
        $defIndex = $i
        break
    }
}

if ($defIndex -gt 1) {
    $termIndex = $defIndex - 1
    $startPos = $paras.Item($termIndex).Range.Start
    $endPos = $paras.Item($defIndex).Range.End
    $d.Range($startPos, $endPos).Delete()
}

Write-Host "found1=$found1 found2=$found2 found3=$found3 removedEntryAt=$defIndex"
